$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "(in per cent)" -> "(in percent)" label in C2
$ws.Range("C2").Value = "(in percent)"

# 2. Add a new year column T (2023) mirroring the existing S column (2022)
#    Header year
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 2023

#    Renewable energy share (%) - also update the 2022 figure
$ws.Range("S5").Value = 29.9
$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("T5").Value = 29.5

#    Hydropower electricity production
$ws.Range("S6").Copy($ws.Range("T6"))
$ws.Range("T6").Value = 12030.6

# 3. Normalise the column widths for the year columns (D:T) to a uniform width
$ws.Range("D1:T1").ColumnWidth = 8.166666666666666

# 4. Reset the sheet selection back to the default top-left cell
$ws.Range("A1").Select() | Out-Null
